$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "H2O"
$ws.Range("B1").Value = "H2"
$ws.Range("C1").Value = "O2"
$ws.Range("D1").ClearContents()

$ws.Range("A2").Select()
